$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("room")
$ws.Range("C2").Value = "adfd"
